$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_WVR = $wb.Worksheets.Item("WVR")

# ALC row 33
$ws_ALC.Cells.Item(33, 8).Value = 752
$ws_ALC.Cells.Item(33, 9).Value = 752
$ws_ALC.Cells.Item(33, 11).Value = 752
$ws_ALC.Cells.Item(33, 13).Value = -523

# ALC row 53
$ws_ALC.Cells.Item(53, 8).Value = 319
$ws_ALC.Cells.Item(53, 9).Value = 301.16666
$ws_ALC.Cells.Item(53, 10).Value = 426
$ws_ALC.Cells.Item(53, 11).Value = 301.16666
$ws_ALC.Cells.Item(53, 12).Value = 426
$ws_ALC.Cells.Item(53, 13).Value = 335.83334
$ws_ALC.Cells.Item(53, 14).Value = -1700

# ALC row 132
$ws_ALC.Cells.Item(132, 8).Value = 3292.926
$ws_ALC.Cells.Item(132, 9).Value = 2454.7917
$ws_ALC.Cells.Item(132, 11).Value = 7364.375100000001
$ws_ALC.Cells.Item(132, 13).Value = -4834.375100000001

# ALC row 133
$ws_ALC.Cells.Item(133, 8).Value = 99995
$ws_ALC.Cells.Item(133, 10).Value = 99995
$ws_ALC.Cells.Item(133, 12).Value = 99995
$ws_ALC.Cells.Item(133, 14).Value = -110115

# ARM row 74
$ws_ARM.Cells.Item(74, 8).Value = 6999.8335
$ws_ARM.Cells.Item(74, 9).Value = 3999.5
$ws_ARM.Cells.Item(74, 11).Value = 3999.5
$ws_ARM.Cells.Item(74, 13).Value = -3125.5

# ARM row 77
$ws_ARM.Cells.Item(77, 8).Value = 6999.8335
$ws_ARM.Cells.Item(77, 9).Value = 3999.5
$ws_ARM.Cells.Item(77, 11).Value = 19997.5
$ws_ARM.Cells.Item(77, 13).Value = -15629.5

# ARM row 93
$ws_ARM.Cells.Item(93, 8).Value = 0
$ws_ARM.Cells.Item(93, 10).Value = 0
$ws_ARM.Cells.Item(93, 12).Value = 0
$ws_ARM.Cells.Item(93, 14).ClearContents()

# ARM row 97
$ws_ARM.Cells.Item(97, 8).Value = 2389.7693
$ws_ARM.Cells.Item(97, 9).Value = 1470.875
$ws_ARM.Cells.Item(97, 11).Value = 1470.875
$ws_ARM.Cells.Item(97, 13).Value = -974.875

# ARM row 122
$ws_ARM.Cells.Item(122, 8).Value = 2485.2856
$ws_ARM.Cells.Item(122, 9).Value = 2324.25
$ws_ARM.Cells.Item(122, 10).Value = 2700
$ws_ARM.Cells.Item(122, 11).Value = 6972.75
$ws_ARM.Cells.Item(122, 12).Value = 8100
$ws_ARM.Cells.Item(122, 13).Value = -4522.75
$ws_ARM.Cells.Item(122, 14).Value = -13000

# ARM row 132
$ws_ARM.Cells.Item(132, 8).Value = 3150.1667
$ws_ARM.Cells.Item(132, 9).Value = 2962.3333
$ws_ARM.Cells.Item(132, 10).Value = 3338
$ws_ARM.Cells.Item(132, 11).Value = 8886.999899999999
$ws_ARM.Cells.Item(132, 12).Value = 10014
$ws_ARM.Cells.Item(132, 13).Value = -6356.999899999999
$ws_ARM.Cells.Item(132, 14).Value = -15074

# BSM row 23
$ws_BSM.Cells.Item(23, 8).Value = 1000
$ws_BSM.Cells.Item(23, 9).Value = 0
$ws_BSM.Cells.Item(23, 11).Value = 0
$ws_BSM.Cells.Item(23, 13).ClearContents()

# BSM row 26
$ws_BSM.Cells.Item(26, 8).Value = 10000
$ws_BSM.Cells.Item(26, 9).Value = 10000
$ws_BSM.Cells.Item(26, 10).Value = 0
$ws_BSM.Cells.Item(26, 11).Value = 10000
$ws_BSM.Cells.Item(26, 12).Value = 0
$ws_BSM.Cells.Item(26, 13).Value = -9708
$ws_BSM.Cells.Item(26, 14).ClearContents()

# BSM row 99
$ws_BSM.Cells.Item(99, 8).Value = 0
$ws_BSM.Cells.Item(99, 9).Value = 0
$ws_BSM.Cells.Item(99, 11).Value = 0
$ws_BSM.Cells.Item(99, 13).ClearContents()

# CRP row 4
$ws_CRP.Cells.Item(4, 8).Value = 4812
$ws_CRP.Cells.Item(4, 9).Value = 0
$ws_CRP.Cells.Item(4, 10).Value = 4812
$ws_CRP.Cells.Item(4, 11).Value = 0
$ws_CRP.Cells.Item(4, 12).Value = 4812
$ws_CRP.Cells.Item(4, 13).ClearContents()
$ws_CRP.Cells.Item(4, 14).Value = -5036

# CRP row 111
$ws_CRP.Cells.Item(111, 8).Value = 99995
$ws_CRP.Cells.Item(111, 10).Value = 99995
$ws_CRP.Cells.Item(111, 12).Value = 99995
$ws_CRP.Cells.Item(111, 14).Value = -108175

# CRP row 132
$ws_CRP.Cells.Item(132, 8).Value = 4933
$ws_CRP.Cells.Item(132, 9).Value = 4933
$ws_CRP.Cells.Item(132, 11).Value = 14799
$ws_CRP.Cells.Item(132, 13).Value = -12269

# CRP row 138
$ws_CRP.Cells.Item(138, 8).Value = 99987
$ws_CRP.Cells.Item(138, 10).Value = 99987
$ws_CRP.Cells.Item(138, 12).Value = 99987
$ws_CRP.Cells.Item(138, 14).Value = -110267

# CUL row 4
$ws_CUL.Cells.Item(4, 8).Value = 375518.56
$ws_CUL.Cells.Item(4, 9).Value = 357664.16
$ws_CUL.Cells.Item(4, 11).Value = 1072992.48
$ws_CUL.Cells.Item(4, 13).Value = -1072880.48

# CUL row 7
$ws_CUL.Cells.Item(7, 8).Value = 0
$ws_CUL.Cells.Item(7, 9).Value = 0
$ws_CUL.Cells.Item(7, 11).Value = 0
$ws_CUL.Cells.Item(7, 13).ClearContents()

# CUL row 9
$ws_CUL.Cells.Item(9, 8).Value = 431.66666
$ws_CUL.Cells.Item(9, 9).Value = 431.66666
$ws_CUL.Cells.Item(9, 11).Value = 1294.99998
$ws_CUL.Cells.Item(9, 13).Value = -1070.99998

# CUL row 16
$ws_CUL.Cells.Item(16, 8).Value = 137
$ws_CUL.Cells.Item(16, 9).Value = 137
$ws_CUL.Cells.Item(16, 10).Value = 0
$ws_CUL.Cells.Item(16, 11).Value = 411
$ws_CUL.Cells.Item(16, 12).Value = 0
$ws_CUL.Cells.Item(16, 13).Value = -238
$ws_CUL.Cells.Item(16, 14).ClearContents()

# CUL row 55
$ws_CUL.Cells.Item(55, 8).Value = 1666.6666
$ws_CUL.Cells.Item(55, 9).Value = 1400
$ws_CUL.Cells.Item(55, 10).Value = 3000
$ws_CUL.Cells.Item(55, 11).Value = 4200
$ws_CUL.Cells.Item(55, 12).Value = 9000
$ws_CUL.Cells.Item(55, 13).Value = -4023
$ws_CUL.Cells.Item(55, 14).Value = -9354

# CUL row 93
$ws_CUL.Cells.Item(93, 8).Value = 2000
$ws_CUL.Cells.Item(93, 10).Value = 0
$ws_CUL.Cells.Item(93, 12).Value = 0
$ws_CUL.Cells.Item(93, 14).ClearContents()

# CUL row 98
$ws_CUL.Cells.Item(98, 8).Value = 1960
$ws_CUL.Cells.Item(98, 10).Value = 1996
$ws_CUL.Cells.Item(98, 12).Value = 5988
$ws_CUL.Cells.Item(98, 14).Value = -8984

# CUL row 131
$ws_CUL.Cells.Item(131, 8).Value = 2175.5
$ws_CUL.Cells.Item(131, 10).Value = 3218.75
$ws_CUL.Cells.Item(131, 12).Value = 9656.25
$ws_CUL.Cells.Item(131, 14).Value = -19736.25

# GSM row 97
$ws_GSM.Cells.Item(97, 8).Value = 2368.75
$ws_GSM.Cells.Item(97, 9).Value = 2368.75
$ws_GSM.Cells.Item(97, 11).Value = 2368.75
$ws_GSM.Cells.Item(97, 13).Value = -1872.75

# LTW row 13
$ws_LTW.Cells.Item(13, 8).Value = 7254.6
$ws_LTW.Cells.Item(13, 9).Value = 7254.6
$ws_LTW.Cells.Item(13, 11).Value = 7254.6
$ws_LTW.Cells.Item(13, 13).Value = -7114.6

# LTW row 16
$ws_LTW.Cells.Item(16, 8).Value = 1133.3334
$ws_LTW.Cells.Item(16, 9).Value = 1133.3334
$ws_LTW.Cells.Item(16, 11).Value = 1133.3334
$ws_LTW.Cells.Item(16, 13).Value = -963.3334

# LTW row 46
$ws_LTW.Cells.Item(46, 8).Value = 4626.548
$ws_LTW.Cells.Item(46, 9).Value = 3249.75
$ws_LTW.Cells.Item(46, 10).Value = 4771.4736
$ws_LTW.Cells.Item(46, 11).Value = 3249.75
$ws_LTW.Cells.Item(46, 12).Value = 4771.4736
$ws_LTW.Cells.Item(46, 13).Value = -3061.75
$ws_LTW.Cells.Item(46, 14).Value = -5147.4736

# LTW row 82
$ws_LTW.Cells.Item(82, 8).Value = 1983.625
$ws_LTW.Cells.Item(82, 9).Value = 2190
$ws_LTW.Cells.Item(82, 11).Value = 2190
$ws_LTW.Cells.Item(82, 13).Value = -1829

# LTW row 85
$ws_LTW.Cells.Item(85, 8).Value = 1983.625
$ws_LTW.Cells.Item(85, 9).Value = 2190
$ws_LTW.Cells.Item(85, 11).Value = 2190
$ws_LTW.Cells.Item(85, 13).Value = -942

# LTW row 100
$ws_LTW.Cells.Item(100, 8).Value = 2861.1
$ws_LTW.Cells.Item(100, 9).Value = 2956.7778
$ws_LTW.Cells.Item(100, 11).Value = 2956.7778
$ws_LTW.Cells.Item(100, 13).Value = -2415.7778

# LTW row 122
$ws_LTW.Cells.Item(122, 8).Value = 8950
$ws_LTW.Cells.Item(122, 9).Value = 8600
$ws_LTW.Cells.Item(122, 11).Value = 25800
$ws_LTW.Cells.Item(122, 13).Value = -23350

# WVR row 104
$ws_WVR.Cells.Item(104, 8).Value = 33998.5
$ws_WVR.Cells.Item(104, 10).Value = 33998.5
$ws_WVR.Cells.Item(104, 12).Value = 33998.5
$ws_WVR.Cells.Item(104, 14).Value = -40986.5

# WVR row 107
$ws_WVR.Cells.Item(107, 8).Value = 2200.1428
$ws_WVR.Cells.Item(107, 9).Value = 1149.75
$ws_WVR.Cells.Item(107, 10).Value = 3600.6667
$ws_WVR.Cells.Item(107, 11).Value = 3449.25
$ws_WVR.Cells.Item(107, 12).Value = 10802.0001
$ws_WVR.Cells.Item(107, 13).Value = -1529.25
$ws_WVR.Cells.Item(107, 14).Value = -14642.0001

